# Update DF_Partial_Result sheet with the re-run GridSearchCV results.
# Rows are now sorted by mean_test_score (desc) from the new run, which changes
# the param-index (A), score columns (C/D), and rank (E) for every data row,
# and the associated params text shown in column B.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$results = @(
    @(29, "{'anOptimizer': 'adam', 'batch_size': 20, 'epochs': 10, 'hidUnit': 256, 'outActivation': 'sigmoid'}", 0.9139906565348307, 0.003887575425954086, 1),
    @(21, "{'anOptimizer': 'adam', 'batch_size': 25, 'epochs': 10, 'hidUnit': 256, 'outActivation': 'sigmoid'}", 0.908005436261495, 0.002555320589502078, 2),
    @(30, "{'anOptimizer': 'adam', 'batch_size': 20, 'epochs': 10, 'hidUnit': 128, 'outActivation': 'softmax'}", 0.9079903761545817, 0.003902322145350392, 3),
    @(28, "{'anOptimizer': 'adam', 'batch_size': 20, 'epochs': 10, 'hidUnit': 256, 'outActivation': 'softmax'}", 0.9059550960858663, 0.008008144880462644, 4),
    @(12, "{'anOptimizer': 'rmsprop', 'batch_size': 20, 'epochs': 10, 'hidUnit': 256, 'outActivation': 'softmax'}", 0.9039650162061056, 0.008827908395781417, 5),
    @(16, "{'anOptimizer': 'adam', 'batch_size': 25, 'epochs': 5, 'hidUnit': 256, 'outActivation': 'softmax'}", 0.9029247562090555, 0.01091472741616778, 6),
    @(14, "{'anOptimizer': 'rmsprop', 'batch_size': 20, 'epochs': 10, 'hidUnit': 128, 'outActivation': 'softmax'}", 0.901025136311849, 0.008132732355088356, 7),
    @(23, "{'anOptimizer': 'adam', 'batch_size': 25, 'epochs': 10, 'hidUnit': 128, 'outActivation': 'sigmoid'}", 0.8989747961362203, 0.006036515920283998, 8),
    @(31, "{'anOptimizer': 'adam', 'batch_size': 20, 'epochs': 10, 'hidUnit': 128, 'outActivation': 'sigmoid'}", 0.8989446957906088, 0.01070936497005477, 9),
    @(15, "{'anOptimizer': 'rmsprop', 'batch_size': 20, 'epochs': 10, 'hidUnit': 128, 'outActivation': 'sigmoid'}", 0.8978591958681742, 0.02273429654968474, 10),
    @(25, "{'anOptimizer': 'adam', 'batch_size': 20, 'epochs': 5, 'hidUnit': 256, 'outActivation': 'sigmoid'}", 0.8928840359052023, 0.01795231191475042, 11),
    @(13, "{'anOptimizer': 'rmsprop', 'batch_size': 20, 'epochs': 10, 'hidUnit': 256, 'outActivation': 'sigmoid'}", 0.8908035357793173, 0.02872724500098014, 12),
    @(22, "{'anOptimizer': 'adam', 'batch_size': 25, 'epochs': 10, 'hidUnit': 128, 'outActivation': 'softmax'}", 0.890079915523529, 0.01283822570121066, 13),
    @(19, "{'anOptimizer': 'adam', 'batch_size': 25, 'epochs': 5, 'hidUnit': 128, 'outActivation': 'sigmoid'}", 0.8879541754722595, 0.008110013959894309, 14),
    @(26, "{'anOptimizer': 'adam', 'batch_size': 20, 'epochs': 5, 'hidUnit': 128, 'outActivation': 'softmax'}", 0.8859640757242838, 0.01547915948105867, 15),
    @(24, "{'anOptimizer': 'adam', 'batch_size': 20, 'epochs': 5, 'hidUnit': 256, 'outActivation': 'softmax'}", 0.8799789150555929, 0.004719762923320736, 16),
    @(27, "{'anOptimizer': 'adam', 'batch_size': 20, 'epochs': 5, 'hidUnit': 128, 'outActivation': 'sigmoid'}", 0.8790742953618368, 0.01079003853631027, 17),
    @(5, "{'anOptimizer': 'rmsprop', 'batch_size': 25, 'epochs': 10, 'hidUnit': 256, 'outActivation': 'sigmoid'}", 0.8711744546890259, 0.04822146359487307, 18),
    @(20, "{'anOptimizer': 'adam', 'batch_size': 25, 'epochs': 10, 'hidUnit': 256, 'outActivation': 'softmax'}", 0.8711141347885132, 0.01771889732193794, 19),
    @(6, "{'anOptimizer': 'rmsprop', 'batch_size': 25, 'epochs': 10, 'hidUnit': 128, 'outActivation': 'softmax'}", 0.8690034747123718, 0.002486157563074196, 20),
    @(18, "{'anOptimizer': 'adam', 'batch_size': 25, 'epochs': 5, 'hidUnit': 128, 'outActivation': 'softmax'}", 0.8680687546730042, 0.01858779946813298, 21),
    @(17, "{'anOptimizer': 'adam', 'batch_size': 25, 'epochs': 5, 'hidUnit': 256, 'outActivation': 'sigmoid'}", 0.8640886346499125, 0.02421498586370427, 22),
    @(7, "{'anOptimizer': 'rmsprop', 'batch_size': 25, 'epochs': 10, 'hidUnit': 128, 'outActivation': 'sigmoid'}", 0.8638022343317667, 0.02962320974625736, 23),
    @(10, "{'anOptimizer': 'rmsprop', 'batch_size': 20, 'epochs': 5, 'hidUnit': 128, 'outActivation': 'softmax'}", 0.8631237943967184, 0.02806388880921575, 24),
    @(1, "{'anOptimizer': 'rmsprop', 'batch_size': 25, 'epochs': 5, 'hidUnit': 256, 'outActivation': 'sigmoid'}", 0.8610884745915731, 0.03979730990531141, 25),
    @(9, "{'anOptimizer': 'rmsprop', 'batch_size': 20, 'epochs': 5, 'hidUnit': 256, 'outActivation': 'sigmoid'}", 0.8599427342414856, 0.01267626852784687, 26),
    @(11, "{'anOptimizer': 'rmsprop', 'batch_size': 20, 'epochs': 5, 'hidUnit': 128, 'outActivation': 'sigmoid'}", 0.8581185340881348, 0.02146008484331551, 27),
    @(0, "{'anOptimizer': 'rmsprop', 'batch_size': 25, 'epochs': 5, 'hidUnit': 256, 'outActivation': 'softmax'}", 0.8408111333847046, 0.05346210040609972, 28),
    @(4, "{'anOptimizer': 'rmsprop', 'batch_size': 25, 'epochs': 10, 'hidUnit': 256, 'outActivation': 'softmax'}", 0.8141263524691263, 0.08834974926229308, 29),
    @(8, "{'anOptimizer': 'rmsprop', 'batch_size': 20, 'epochs': 5, 'hidUnit': 256, 'outActivation': 'softmax'}", 0.8067993521690369, 0.02879368857426483, 30),
    @(2, "{'anOptimizer': 'rmsprop', 'batch_size': 25, 'epochs': 5, 'hidUnit': 128, 'outActivation': 'softmax'}", 0.8039047122001648, 0.01900706778250232, 31),
    @(3, "{'anOptimizer': 'rmsprop', 'batch_size': 25, 'epochs': 5, 'hidUnit': 128, 'outActivation': 'sigmoid'}", 0.7999396920204163, 0.008528330502482386, 32)
)

for ($i = 0; $i -lt $results.Length; $i++) {
    $row = $i + 2
    $rec = $results[$i]
    $ws.Cells.Item($row, 1).Value2 = $rec[0]
    $ws.Cells.Item($row, 2).Value2 = $rec[1]
    $ws.Cells.Item($row, 3).Value2 = $rec[2]
    $ws.Cells.Item($row, 4).Value2 = $rec[3]
    $ws.Cells.Item($row, 5).Value2 = $rec[4]
}
